$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 392, shifting existing rows 392:462 down to 393:463
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with the new weekly record
$ws.Cells.Item(392, 1).Value = 8
$ws.Cells.Item(392, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(392, 3).Value = "Coquimbo"
$ws.Cells.Item(392, 4).Value = 44785
$ws.Cells.Item(392, 5).Value = 4
$ws.Cells.Item(392, 6).Value = 100112009
$ws.Cells.Item(392, 7).Value = "Acelga"
$ws.Cells.Item(392, 8).Value = "Sin especificar"
$ws.Cells.Item(392, 9).Value = "Segunda"
$ws.Cells.Item(392, 10).Value = 1460
$ws.Cells.Item(392, 11).Value = 600
$ws.Cells.Item(392, 12).Value = 650
$ws.Cells.Item(392, 13).Value = 625
$ws.Cells.Item(392, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(392, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(392, 16).Value = 312
$ws.Cells.Item(392, 17).Value = 2
$ws.Cells.Item(392, 18).Value = "Hortaliza"
